{"js": "const body = context.document.body;\n\n// 1. Insert the new sentence about Thomas's alignment-functionality work into the\n//    first \"Thomas started by working...\" paragraph, right after the splash-screen\n//    sentence and before \"Other duties included:\".\nconst target = body.search(\n  \"a how to play splash screen.  Other duties included:\",\n  { matchCase: true }\n);\ntarget.load(\"text\");\nawait context.sync();\n\ntarget.items[0].insertText(\n  \"a how to play splash screen. He also designed and implemented an alignment \" +\n    \"functionality for the toolbar but this functionality was not used in the \" +\n    \"final product since it was brought to the attention of the group on the \" +\n    \"last day before the code-freeze. Other duties included:\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 2. The lone \"_GoBack\" bookmark moves from the second \"Finally, Thomas...\"\n//    paragraph (right after \"for the toolbar\") up into the sentence we just\n//    inserted, wrapping \"...before the code-freeze\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst bookmarkAnchor = body.search(\n  \"on the last day before the code-freeze\",\n  { matchCase: true }\n);\nbookmarkAnchor.load(\"text\");\nawait context.sync();\n\nbookmarkAnchor.items[0].insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Insert the new sentence about Thomas's alignment-functionality work into the\n#    first \"Thomas started by working...\" paragraph, right after the splash-screen\n#    sentence and before \"Other duties included:\".\n$r = $d.Content\n$null = $r.Find.Execute(\"a how to play splash screen.  Other duties included:\")\n$r.Text = \"a how to play splash screen. He also designed and implemented an alignment functionality for the toolbar but this functionality was not used in the final product since it was brought to the attention of the group on the last day before the code-freeze. Other duties included:\"\n\n# 2. The lone \"_GoBack\" bookmark moves from the second \"Finally, Thomas...\"\n#    paragraph (right after \"for the toolbar\") up into the sentence we just\n#    inserted, wrapping up right after \"...before the code-freeze\".\n$d.Bookmarks(\"_GoBack\").Delete()\n\n$r2 = $d.Content\n$null = $r2.Find.Execute(\"on the last day before the code-freeze\")\n$r2.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $r2)\n"}
